# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (F column) figures across the sheets.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 3263
$ws1.Range("F4").Value  = 2014
$ws1.Range("F6").Value  = 114
$ws1.Range("F7").Value  = 3097
$ws1.Range("F9").Value  = 299
$ws1.Range("F10").Value = 42
$ws1.Range("F15").Value = 10217
$ws1.Range("F20").Value = 8115
$ws1.Range("F21").Value = 12724
$ws1.Range("F33").Value = 8015
$ws1.Range("F34").Value = 1618
$ws1.Range("F38").Value = 4629
$ws1.Range("F39").Value = 1471
$ws1.Range("F40").Value = 74
$ws1.Range("F41").Value = 386

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 1202
$ws2.Range("F12").Value = 32
$ws2.Range("F15").Value = 15
$ws2.Range("F18").Value = 22

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 647
$ws3.Range("F4").Value = 229
$ws3.Range("F5").Value = 23

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 647
$ws4.Range("F4").Value  = 3263
$ws4.Range("F6").Value  = 2014
$ws4.Range("F9").Value  = 23
$ws4.Range("F10").Value = 3097
$ws4.Range("F13").Value = 42
$ws4.Range("F18").Value = 10217
$ws4.Range("F22").Value = 8116
$ws4.Range("F23").Value = 12724
$ws4.Range("F33").Value = 32
$ws4.Range("F36").Value = 8015
$ws4.Range("F40").Value = 4629
$ws4.Range("F44").Value = 22
